$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1159.3636
$ws.Range("I70").Value = 1292
$ws.Range("J70").Value = 1000.2
$ws.Range("K70").Value = 3876
$ws.Range("L70").Value = 3000.6
$ws.Range("M70").Value = -3606
$ws.Range("N70").Value = -3540.6
$ws.Range("H73").Value = 1159.3636
$ws.Range("I73").Value = 1292
$ws.Range("J73").Value = 1000.2
$ws.Range("K73").Value = 3876
$ws.Range("L73").Value = 3000.6
$ws.Range("M73").Value = -2940
$ws.Range("N73").Value = -4872.6
$ws.Range("H74").Value = 5072.727
$ws.Range("I74").Value = 4780
$ws.Range("J74").Value = 8000
$ws.Range("K74").Value = 4780
$ws.Range("L74").Value = 8000
$ws.Range("M74").Value = -3844
$ws.Range("N74").Value = -9872
$ws.Range("H76").Value = 3349.2666
$ws.Range("I76").Value = 3125.3635
$ws.Range("J76").Value = 3965
$ws.Range("K76").Value = 3125.3635
$ws.Range("L76").Value = 3965
$ws.Range("M76").Value = -2810.3635
$ws.Range("N76").Value = -4595
$ws.Range("H77").Value = 5072.727
$ws.Range("I77").Value = 4780
$ws.Range("J77").Value = 8000
$ws.Range("K77").Value = 23900
$ws.Range("L77").Value = 40000
$ws.Range("M77").Value = -19220
$ws.Range("N77").Value = -49360
$ws.Range("H79").Value = 3349.2666
$ws.Range("I79").Value = 3125.3635
$ws.Range("J79").Value = 3965
$ws.Range("K79").Value = 3125.3635
$ws.Range("L79").Value = 3965
$ws.Range("M79").Value = -2033.3635
$ws.Range("N79").Value = -6149
$ws.Range("H80").Value = 96170.766
$ws.Range("I80").Value = 143614.58
$ws.Range("K80").Value = 430843.74
$ws.Range("M80").Value = -429845.74
$ws.Range("H83").Value = 96170.766
$ws.Range("I83").Value = 143614.58
$ws.Range("K83").Value = 1292531.22
$ws.Range("M83").Value = -1287539.22
$ws.Range("H138").Value = 4712
$ws.Range("I138").Value = 2029.25
$ws.Range("J138").Value = 6053.375
$ws.Range("K138").Value = 6087.75
$ws.Range("L138").Value = 18160.125
$ws.Range("M138").Value = -947.75
$ws.Range("N138").Value = -28440.125

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2637.05
$ws.Range("I61").Value = 1648.875
$ws.Range("J61").Value = 3295.8333
$ws.Range("K61").Value = 1648.875
$ws.Range("L61").Value = 3295.8333
$ws.Range("M61").Value = -1436.875
$ws.Range("N61").Value = -3719.8333
$ws.Range("H88").Value = 2670.3
$ws.Range("I88").Value = 2689.2222
$ws.Range("K88").Value = 2689.2222
$ws.Range("M88").Value = -2283.2222
$ws.Range("H91").Value = 2670.3
$ws.Range("I91").Value = 2689.2222
$ws.Range("K91").Value = 2689.2222
$ws.Range("M91").Value = -1285.2222
$ws.Range("H109").Value = 33350
$ws.Range("J109").Value = 33350
$ws.Range("L109").Value = 33350
$ws.Range("N109").Value = -36124
$ws.Range("H136").Value = 2637.05
$ws.Range("I136").Value = 1648.875
$ws.Range("J136").Value = 3295.8333
$ws.Range("K136").Value = 4946.625
$ws.Range("L136").Value = 9887.499899999999
$ws.Range("M136").Value = -2396.625
$ws.Range("N136").Value = -14987.4999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 47330
$ws.Range("J62").Value = 46995
$ws.Range("L62").Value = 46995
$ws.Range("N62").Value = -48367
$ws.Range("H65").Value = 47330
$ws.Range("J65").Value = 46995
$ws.Range("L65").Value = 140985
$ws.Range("N65").Value = -147849

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 27330
$ws.Range("J87").Value = 27330
$ws.Range("L87").Value = 27330
$ws.Range("N87").Value = -29702
$ws.Range("H90").Value = 27330
$ws.Range("J90").Value = 27330
$ws.Range("L90").Value = 81990
$ws.Range("N90").Value = -93846
$ws.Range("H108").Value = 29163
$ws.Range("J108").Value = 29163
$ws.Range("L108").Value = 29163
$ws.Range("N108").Value = -36843

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 85
$ws.Range("J2").Value = 86.666664
$ws.Range("L2").Value = 519.999984
$ws.Range("N2").Value = -745.999984
$ws.Range("H12").Value = 93.875
$ws.Range("J12").Value = 93.875
$ws.Range("L12").Value = 281.625
$ws.Range("N12").Value = -627.625
$ws.Range("H17").Value = 1412.75
$ws.Range("I17").Value = 1200
$ws.Range("J17").Value = 2051
$ws.Range("K17").Value = 3600
$ws.Range("L17").Value = 6153
$ws.Range("M17").Value = -3431
$ws.Range("N17").Value = -6491
$ws.Range("H20").Value = 999.86664
$ws.Range("I20").Value = 999.86664
$ws.Range("K20").Value = 2999.59992
$ws.Range("M20").Value = -2772.59992
$ws.Range("H34").Value = 2166.5557
$ws.Range("J34").Value = 2742.7144
$ws.Range("L34").Value = 8228.143199999999
$ws.Range("N34").Value = -8396.143199999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2691.6667
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 2830
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 2830
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4826
$ws.Range("H83").Value = 2691.6667
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 2830
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 14150
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -24134

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H68").Value = 4650.25
$ws.Range("I68").Value = 1500.5
$ws.Range("J68").Value = 7800
$ws.Range("K68").Value = 1500.5
$ws.Range("L68").Value = 7800
$ws.Range("M68").Value = -751.5
$ws.Range("N68").Value = -9298
$ws.Range("H71").Value = 4650.25
$ws.Range("I71").Value = 1500.5
$ws.Range("J71").Value = 7800
$ws.Range("K71").Value = 7502.5
$ws.Range("L71").Value = 39000
$ws.Range("M71").Value = -3758.5
$ws.Range("N71").Value = -46488
$ws.Range("H82").Value = 1495
$ws.Range("J82").Value = 1368
$ws.Range("L82").Value = 1368
$ws.Range("N82").Value = -2090
$ws.Range("H85").Value = 1495
$ws.Range("J85").Value = 1368
$ws.Range("L85").Value = 1368
$ws.Range("N85").Value = -3864
$ws.Range("H136").Value = 1665.8
$ws.Range("I136").Value = 1544.7693
$ws.Range("J136").Value = 2452.5
$ws.Range("K136").Value = 4634.3079
$ws.Range("L136").Value = 7357.5
$ws.Range("M136").Value = -2084.3079
$ws.Range("N136").Value = -12457.5
